$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Latest_stocks sheet (sheet1) - quantity updates and highlight fills
# ---------------------------------------------------------------------
$wsStocks = $wb.Worksheets.Item("Latest_stocks")

# Quantity corrections
$wsStocks.Range("C3").Value = 9
$wsStocks.Range("C4").Value = 10
$wsStocks.Range("C13").Value = 0
$wsStocks.Range("C15").Value = 5
$wsStocks.Range("C16").Value = 6
$wsStocks.Range("C17").Value = 0
$wsStocks.Range("C20").Value = 0
$wsStocks.Range("C28").Value = 66
$wsStocks.Range("C31").Value = 62
$wsStocks.Range("C33").Value = 8

# Rename "150mm ventilation fan" entry to "150mm exhaust fan"
$wsStocks.Range("B24").Value = "150mm exhaust fan"

# Highlight the TPW product rows (B14:B20) in yellow
$wsStocks.Range("B14").Interior.Color = 65535
$wsStocks.Range("B15").Interior.Color = 65535
$wsStocks.Range("B16").Interior.Color = 65535
$wsStocks.Range("B17").Interior.Color = 65535
$wsStocks.Range("B18").Interior.Color = 65535
$wsStocks.Range("B19").Interior.Color = 65535
$wsStocks.Range("B20").Interior.Color = 65535

# Highlight the Exhaust/Ventilation fan rows (B21:B26) in green
$wsStocks.Range("B21").Interior.Color = 5296274
$wsStocks.Range("B22").Interior.Color = 5296274
$wsStocks.Range("B23").Interior.Color = 5296274
$wsStocks.Range("B24").Interior.Color = 5296274
$wsStocks.Range("B25").Interior.Color = 5296274
$wsStocks.Range("B26").Interior.Color = 5296274

# Shrink print scale
$wsStocks.PageSetup.Zoom = 70

# ---------------------------------------------------------------------
# Sales sheet (sheet4) - amount corrections
# ---------------------------------------------------------------------
$wsSales = $wb.Worksheets.Item("Sales")
$wsSales.Range("F49").Value = 10000
$wsSales.Range("E64").Value = 23880

# ---------------------------------------------------------------------
# Restore the active-window selection/scroll state for each sheet
# ---------------------------------------------------------------------
$wsPrices = $wb.Worksheets.Item("Price_details")

$wsPrices.Activate()
$wsPrices.Range("A10").Select()

$wsSales.Activate()
$wsSales.Range("B71").Select()

$wsStocks.Activate()
$wsStocks.Range("B40").Select()
